# Adding breathcounter and shared_table example forms
#
# In the "survey" sheet, a new "geopoint" example question ("Capture your
# location") is inserted right after the existing "barcode" question
# (currently row 9, just before the "picture" question). Inserting the row
# natively shifts every subsequent row down by one and lets Excel manage the
# shared-string table / row formatting bookkeeping on its own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$ws.Rows.Item(9).Insert()

$ws.Cells.Item(9, 1).Value = "geopoint"
$ws.Cells.Item(9, 3).Value = "geopoint"
$ws.Cells.Item(9, 4).Value = "Capture your location"
